# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (used by the slide master -> the "live" design)
#     Integral  ->  Office Theme
#   ppt/theme/theme2.xml  (used by the notes master)
#     Office Theme  ->  Integral
#
# The only real content difference between the two theme parts is the
# <a:clrScheme> color values (font + format schemes are already identical).
# The PowerPoint object model exposes the presentation's live color scheme
# through Slide.ThemeColorScheme (backed by theme1.xml, the slide-master
# theme), so we repaint it from "Integral" to the stock "Office" palette,
# matching the new theme1.xml contents from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office theme colors, encoded as OLE RGB ints (0x00BBGGRR) the way
# ColorFormat.RGB expects them. Index order follows the standard
# PowerPoint ColorScheme ordinal: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs.Colors(1).RGB = 0x000000    # dk1    -> 000000 (unchanged)
$tcs.Colors(2).RGB = 0xFFFFFF    # lt1    -> FFFFFF (unchanged)
$tcs.Colors(3).RGB = 0x6A5444    # dk2    -> 44546A
$tcs.Colors(4).RGB = 0xE6E6E7    # lt2    -> E7E6E6
$tcs.Colors(5).RGB = 0xD59B5B    # accent1-> 5B9BD5
$tcs.Colors(6).RGB = 0x317DED    # accent2-> ED7D31
$tcs.Colors(7).RGB = 0xA5A5A5    # accent3-> A5A5A5
$tcs.Colors(8).RGB = 0x00C0FF    # accent4-> FFC000
$tcs.Colors(9).RGB = 0xC47244    # accent5-> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6-> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink  -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
